# Trade #48 closed at 2026-02-17 08:39:22 - unknown UNKNOWN +0.000%
# Update the rolled-up stats on Summary / Strategy Status and append the
# newly-closed trade row to both "All Trades" and "MarketMaking".

$wb = $excel.ActiveWorkbook

# --- Summary sheet -------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.57   # Current Capital
$summary.Range("B4").Value = -0.43     # Total P&L $
$summary.Range("B5").Value = -0.18     # Total P&L %
$summary.Range("B6").Value = 48        # Total Trades
$summary.Range("B7").Value = 18        # Winning Trades
$summary.Range("B9").Value = 37.5      # Win Rate %

# --- Strategy Status sheet (MarketMaking row, row 4) ----------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.56999999999999   # Capital
$status.Range("D4").Value = 48                  # Trades
$status.Range("E4").Value = -0.43               # P&L $
$status.Range("F4").Value = -0.43               # P&L %
$status.Range("G4").Value = 37.5                # Win Rate %

# --- Append the new closed trade (#48) to All Trades + MarketMaking -------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A49").Value = 48

    # "2026-02-17" looks like a date, so Excel would otherwise silently
    # coerce it to a date serial number on assignment. Force the cell to
    # text first, then clear the format override again afterwards so the
    # cell keeps the workbook's default (unstyled) appearance, matching
    # how the other inline-string cells in this column are stored.
    $dateCell = $ws.Range("B49")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.ClearFormats()

    $ws.Range("C49").Value = "08:39:16"
    $ws.Range("D49").Value = "MarketMaking"
    $ws.Range("E49").Value = "DOWN"
    $ws.Range("F49").Value = 0.14
    $ws.Range("G49").Value = 0.15
    $ws.Range("H49").Value = "CLOSED"
    $ws.Range("I49").Value = 7.1429
    $ws.Range("J49").Value = 0.01
    $ws.Range("K49").Value = 99.56999999999999
    $ws.Range("L49").Value = 0
    $ws.Range("M49").Value = 0
    $ws.Range("N49").Value = 0.6
    $ws.Range("O49").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P49").Value = "early_exit"
    $ws.Range("Q49").Value = 0.14
}
